$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 15 (brand-new "SD-N" run row) — values entered first, in the same
# left-to-right-ish order the author typed them in (matches shared-string
# table growth order: A,B,E,D,I,G,H,M,F,J,K).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "SD-N"
$ws.Range("B15").Value = "PPO use step distance reward + multiply critic lr + train every episode + normalized reward vs. Random"
$ws.Range("E15").Value = "老大爷遛弯"
$ws.Range("D15").Value = "对着墙跑，不怎么能移动"
$ws.Range("I15").Value = "卡很久之后能过"
$ws.Range("G15").Value = "完全不行"
$ws.Range("H15").Value = "转圈"
$ws.Range("M15").Value = "四处碰壁"
$ws.Range("F15").Value = "转圈，上下碰"
$ws.Range("J15").Value = "不行"
$ws.Range("K15").Value = "稳"
$ws.Range("L15").Value = "转圈"

# Formatting for row 15, copied from rows that already use the right styles
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("M15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 13 (existing "SD" run row): fill in the newly-finished per-map columns
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = "会折返很多，有时不能过"
$ws.Range("E13").Value = "即使快到了也要回去绕"
$ws.Range("D13").Value = "还算顺利"
$ws.Range("G13").Value = "not a chance"
$ws.Range("H13").Value = "受random干扰很大，但有一次居然靠自己过了"
$ws.Range("I13").Value = "random没干扰就能过"
$ws.Range("L13").Value = "干扰很大，但是有时候能走的很远，几乎到终点"
$ws.Range("M13").Value = "一开始卡半天，后来也只能过个弯"

# Match the formatting used by the other cells in that row / column family
$ws.Range("C3").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("I10").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("I10").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("L14").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$ws.Range("M14").Copy()
$ws.Range("M13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# View state: scroll/zoom/selection updated by the author after this edit
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 83
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G10").Select()
